$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) cells to Text format so that numeric-looking
# strings (e.g. "213.19") are not reinterpreted as numbers by Excel, then restore
# the original (default) formatting once the text values are in place.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.500.89"
$ws.Range("E2").Value = "  +4.11%  "
$ws.Range("D3").Value = "1.593.76"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").Value = "213.19"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  +6.67%  "
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "26.70"
$ws.Range("E8").Value = "  +11.96%  "
$ws.Range("D9").Value = "0.249"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").Value = "0.0596"
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("D11").Value = "0.0910"
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").Value = "1.818.48"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "1.575.86"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "29.486.15"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.529"
$ws.Range("E15").Value = "  +3.76%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "3.75"
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("D17").Value = "63.21"
$ws.Range("E17").Value = "  +4.06%  "
$ws.Range("D18").Value = "241.55"
$ws.Range("E18").Value = "  +6.19%  "
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("D21").Value = "0.994"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "4.04"
$ws.Range("E22").Value = "  +3.23%  "
$ws.Range("D23").Value = "9.28"
$ws.Range("E23").Value = "  +4.07%  "
$ws.Range("D25").Value = "154.73"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "15.25"
$ws.Range("E26").Value = "  +3.59%  "
$ws.Range("E27").Value = "  +5.40%  "
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "0.0472"
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "1.428.12"
$ws.Range("E33").Value = "  +3.28%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("D35").Value = "1.05"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("D37").Value = "2.82"
$ws.Range("E37").Value = "  +9.51%  "
$ws.Range("D38").Value = "2.30"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").Value = "0.0166"
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("E40").Value = "  +4.76%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").Value = "54.10"
$ws.Range("E42").Value = "  +28.49%  "
$ws.Range("D43").Value = "0.803"
$ws.Range("E43").Value = "  +3.26%  "
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "0.0471"
$ws.Range("E45").Value = "  +3.48%  "
$ws.Range("D46").Value = "65.18"
$ws.Range("E46").Value = "  +5.24%  "
$ws.Range("D47").Value = "5.36"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "1.729.79"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").Value = "86.37"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "0.834"
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("E51").Value = "  +1.70%  "

$ws.Range("D2:D51").ClearFormats()
